$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Name" column (column B), shifting "Marks" (old column C) left to become column B
$ws.Range("B1").EntireColumn.Delete()

# Update the selection to match the target state
$ws.Range("H10").Select()
